# Price Update: Sat Dec 20 17:12:25 UTC 2025
#
# A new price snapshot column is inserted to the right of the existing
# "current price" column. The former current-price column (B) is copied
# into the new column (C) so it becomes the "previous price" column, and
# column B is refreshed with the latest snapshot (timestamp in row 1,
# prices in the rows below). Only one SKU's price actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Duplicate column B (values + styles) into the new column C. This
#    preserves the header style and every existing price as the
#    "previous" snapshot before we overwrite column B with new data.
$ws.Range("B1:B26").Copy($ws.Range("C1:C26"))

# 2. Match column C's width to column B's (raw stored width of 21).
#    Excel's ColumnWidth property is expressed in character units, which
#    is offset from the raw OOXML width by the default column padding
#    (~0.8333333), so back that out to land exactly on 21.
$ws.Columns.Item(3).ColumnWidth = 21 - 0.8333333333333334

# 3. Refresh column B with the new snapshot.
$ws.Range("B1").Value = "2025-12-20 22:00"

$ws.Range("B19").Value = 1497
